$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B11 text changes from "R40" to the literal text "1".
# Leading apostrophe forces Excel to store it as text (not a number),
# matching the shared-string cell <c r="B11" t="s"> in the target.
$ws.Range("B11").Value = "'1"
